# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds strikeout ("K") counts per game row (rows 2-8).
# Update with the recalculated values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 6
